$d = $word.ActiveDocument

# The last paragraph in the document is an empty bullet item (ListParagraph).
# We need to add three runs of text to it, with a misspelled word
# ("conentratie") wrapped in proofErr spell-check markers, matching the
# exact run-level formatting (rFonts/lang) already used elsewhere in the
# document for this paragraph's run-properties.

$target = $d.Paragraphs.Last
$insertAt = $d.Range($target.Range.Start, $target.Range.Start)

$frag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/><w:lang w:val="nl-BE"/></w:rPr><w:t xml:space="preserve">Alternatief idee: neem bovenaan een CTE waarde van 1, zie dit als een soort van bron die lekt uit de onverzadigde zone die zorgt dat de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/><w:lang w:val="nl-BE"/></w:rPr><w:t>conentratie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/><w:lang w:val="nl-BE"/></w:rPr><w:t xml:space="preserve"> aan de rand altijd dezelfde blijft! </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$insertAt.InsertXML($frag)
